$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 59, duplicating row 58 but with a new ID (EX0058)
$ws.Range("A59").Value = "EX0058"
$ws.Range("B59").Value = $ws.Range("B58").Value
$ws.Range("C59").Value = $ws.Range("C58").Value
$ws.Range("D59").Value = $ws.Range("D58").Value
$ws.Range("E59").Value = $ws.Range("E58").Value
$ws.Range("F59").Value = $ws.Range("F58").Value
$ws.Range("G59").Value = $ws.Range("G58").Value
$ws.Range("H59").Value = $ws.Range("H58").Value
$ws.Range("I59").Value = $ws.Range("I58").Value

# Update the view: scroll so row 37 is at the top, and select A59
$ws.Application.ActiveWindow.ScrollRow = 37
$ws.Range("A59").Select()
